$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the "Periodo Mora" column (E16:E27) from descending
# (2212 .. 2201) to ascending (2201 .. 2212) order, as part of the
# database refresh / "parte 1" update for the new estado de cuenta.
$periods = @("2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# The "Valor Mora" amounts that belonged to period 2201 and 2212 are
# swapped along with the re-sort (they keep following their own period).
$ws.Cells.Item(16, 6).Value = 32000
$ws.Cells.Item(27, 6).Value = 25333
